$wb = $excel.ActiveWorkbook

$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Update 'Recommandations' sheet (sheet1) ---
$dataReco = @(
    @("NEI-CEDA CI", 0, 4, 3750, 950, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES PUBLICS", 0, 8, 3426.56, 111.96, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2391.3, 585.48, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 2158.24, 536.49, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 3, 2105, 705, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE", 0, 2, 1985, 995, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1408.67, 342.68, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1357.82, 338.92, "🟡 Observer", "➖ Neutre"),
    @("ERIUM CI", 0, 1, 810, 810, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 747.8200000000001, 183.58, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 598.08, 150.01, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 587.78, 147.43, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 586.87, 146.76, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                (**)", 0, 2, 538.1799999999999, 269.26, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 491.25, 121.78, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 455.9, 114.81, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                   (**)", 0, 2, 444.8, 222.42, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 387.95, 96.77, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE         (**)", 0, 1, 224.89, 224.89, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE          (**)", 0, 1, 224.67, 224.67, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA BN (BOAB)", 3, 1, 11.94, 3.27, "🟢 Achat", "✅ Renforcer"),
    @("SAFCA CI (SAFC)", 2, 0, 10.87, 7.38, "🟡 Observer", "➖ Neutre"),
    @("ERIUM CI (SIVC)", 1, 0, 7.28, 7.28, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI (SIVC)", 2, 0, 7.09, 7.09, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 1, 0, 6.78, 6.78, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 2, 0, 5.26, 5.26, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 1, 0, 4.97, 4.97, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 2, 1, 4.51, 3.13, "🟡 Observer", "👀 À surveiller"),
    @("BERNABE CI (BNBC)", 1, 1, 3.96, -3.45, "🟡 Observer", "👀 À surveiller"),
    @("SICABLE CI (CABC)", 1, 0, 3, 3, "🟡 Observer", "➖ Neutre"),
    @("SMB CI (SMBC)", 2, 1, 2.98, 2.89, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA ML (BOAM)", 1, 1, 0.83, -2.5, "🟡 Observer", "👀 À surveiller"),
    @("BICI CI (BICC)", 1, 1, 0.33, 2.54, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 0, 1, -1, -1, "🟡 Observer", "➖ Neutre"),
    @("PALM CI (PALC)", 0, 1, -1.54, -1.54, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA SENEGAL (BOAS)", 0, 1, -1.91, -1.91, "🟡 Observer", "➖ Neutre"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 0, 1, -2.06, -2.06, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING SN (TTLS)", 0, 1, -2.35, -2.35, "🟡 Observer", "➖ Neutre"),
    @("ORAGROUP TOGO (ORGT)", 1, 2, -2.4, -3.47, "🟡 Observer", "👀 À surveiller"),
    @("CFAO MOTORS CI (CFAC)", 0, 2, -2.54, -1.14, "🟡 Observer", "👀 À surveiller"),
    @("CORIS BANK INTERNATIONAL (CBIBF)", 0, 1, -2.78, -2.78, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 0, 1, -3.11, -3.11, "🟡 Observer", "➖ Neutre"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -3.45, -3.45, "🟡 Observer", "➖ Neutre"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -4.76, -4.76, "🟡 Observer", "➖ Neutre"),
)

for ($i = 0; $i -lt $dataReco.Count; $i++) {
    $r = $i + 2
    $row = $dataReco[$i]
    $wsReco.Range("A$r").Value = $row[0]
    $wsReco.Range("B$r").Value = $row[1]
    $wsReco.Range("C$r").Value = $row[2]
    $wsReco.Range("D$r").Value = $row[3]
    $wsReco.Range("E$r").Value = $row[4]
    $wsReco.Range("F$r").Value = $row[5]
    $wsReco.Range("G$r").Value = $row[6]
}

# --- Update 'Top_YTD' sheet (sheet2) ---
$dataYtd = @(
    @("BRVM - SERVICES PUBLICS", 10274910.7),
    @("NEI-CEDA CI", 1157773.06),
    @("BRVM - AUTRES SECTEURS", 236951.38),
    @("BRVM - DISTRIBUTION", 167208.1),
    @("AIR LIQUIDE CI", 51420),
    @("BRVM - TRANSPORT", 41689.54),
    @("BRVM - AGRICULTURE", 37194.39),
    @("SUCRIVOIRE", 11835.5),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 6679.69),
    @("BRVM - FINANCES", 3776.28),
)

for ($i = 0; $i -lt $dataYtd.Count; $i++) {
    $r = $i + 2
    $row = $dataYtd[$i]
    $wsYtd.Range("A$r").Value = $row[0]
    $wsYtd.Range("B$r").Value = $row[1]
}

